$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws2 = $wb.Worksheets.Item("Hoja2")

# --- Hoja1: new rows 9-14 ---
$hoja1Rows = @(
  @(70, 1000364),
  @(70, 1000368),
  @(70, 1000370),
  @(70, 1000371),
  @(70, 1000372),
  @(70, 1000373)
)

$r = 9
foreach ($row in $hoja1Rows) {
  $ws1.Range("A$r").Value = "A004"
  $ws1.Range("B$r").Value = $row[0]
  $ws1.Range("C$r").Value = $row[1]
  $r++
}

# Highlight A8 in Hoja1 with yellow fill
$ws1.Range("A8").Interior.Color = 65535

# --- Hoja2: new rows 8-12 ---
$hoja2Rows = @(
  @(100248, 2024, 1000368),
  @(100249, 2024, 1000370),
  @(100250, 2024, 1000371),
  @(100251, 2024, 1000372),
  @(100252, 2024, 1000373)
)

$r = 8
foreach ($row in $hoja2Rows) {
  $ws2.Range("A$r").Value = $row[0]
  $ws2.Range("B$r").Value = $row[1]
  $ws2.Range("C$r").Value = "A004"
  $ws2.Range("D$r").Value = $row[2]
  $r++
}

# Highlight A7 in Hoja2 with yellow fill
$ws2.Range("A7").Interior.Color = 65535

# Set selections to match final state without leaving Hoja1 as the active sheet
$ws2.Activate()
$ws2.Range("A7").Select()
$ws1.Activate()
$ws1.Range("C9").Select()
